# Generate Report for Handback
#
# The localization status report is updated to reflect that the handed-off
# files have now been handed back (and are in sync with en-US):
#   - Status changes from "Ready for handoff" to "Handed back: in sync
#     with en-US" on the Overview sheet and on each language sheet.
#   - Each language sheet (zh-cn, de-de) gets its "Latest Target File"
#     and "Latest Handback File" columns (E, F) populated for the two
#     real source files (a.md / b.md rows), plus a "Latest Handback
#     DateTime" (column G) timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: just the status-text change -------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

function Update-LangSheet {
    param(
        [string]$SheetName,
        [string]$XlfFile,
        [string]$HandbackDateTimeRow2,
        [string]$HandbackDateTimeRow3
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 2 (a.md) ----------------------------------------------------------
    $ws.Range("B2").Value = $newStatus
    $ws.Range("E2").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/a.md", "", "", "a.md") | Out-Null
    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = 15570276

    $ws.Range("F2").Value = $XlfFile
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90df88de8e333fc883e877a686133e687b711794/ol-handback/$XlfFile", "", "", $XlfFile) | Out-Null
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276

    $ws.Range("G2").Value = $HandbackDateTimeRow2
    $ws.Range("H2").Value = "Include"

    # Row 3 (b.md) ------------------------------------------------------------
    $ws.Range("B3").Value = $newStatus
    $ws.Range("E3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/a.md", "", "", "a.md") | Out-Null
    $ws.Range("E3").Font.Underline = 2
    $ws.Range("E3").Font.Color = 15570276

    $ws.Range("F3").Value = $XlfFile
    $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90df88de8e333fc883e877a686133e687b711794/ol-handback/$XlfFile", "", "", $XlfFile) | Out-Null
    $ws.Range("F3").Font.Underline = 2
    $ws.Range("F3").Font.Color = 15570276

    $ws.Range("G3").Value = $HandbackDateTimeRow3
    $ws.Range("H3").Value = "Include"
}

Update-LangSheet "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-02-22 17:31:01" "2016-02-22 17:31:01"
Update-LangSheet "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-02-22 17:31:36" "2016-02-22 17:31:36"
